$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append " Library" to the FacilityName (column A) for each data row (2-22)
for ($r = 2; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    if ($current -ne $null -and $current -notmatch "Library$") {
        $cell.Value2 = "$current Library"
    }
}
